$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sample labels in column A: "sample N" -> "sample_N" (A2:A21)
# (written first so the shared-string table keeps samples before features)
for ($row = 2; $row -le 21; $row++) {
    $sampleNum = $row - 1
    $ws.Cells.Item($row, 1).Value = "sample_$sampleNum"
}

# Rename header feature labels: "feature N" -> "feature_N" (B1:F1)
for ($col = 2; $col -le 6; $col++) {
    $featureNum = $col - 1
    $ws.Cells.Item(1, $col).Value = "feature_$featureNum"
}

# Update selection to B1:F1 with active cell B1
$ws.Range("B1:F1").Select()
